# sszw_vorlage.xlsx update:
#  - column layout rework (A..P) with new widths
#  - row height rework for rows 2-5
#  - "Bauform" header renamed to "Form" (D4) and merged O2:O4 header
#    relabelled "Form" (used to read "keine Uebertrag. an RBC")
#  - new unit row ("m") added for the two "Laenge" columns (G5/H5)
#  - merged cell O2:O3 extended to O2:O4, with O4 taking on the same
#    (bold/centered) formatting as O3
#  - print area extended from A1:O7 to A1:P7
#  - selection moved to J23 (matches the saved view state)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------
# Cell value / label changes
# ---------------------------------------------------------------
$ws.Range("D4").Value = "Form"
$ws.Range("O2").Value = "Form"
$ws.Range("G5").Value = "m"
$ws.Range("H5").Value = "m"

# ---------------------------------------------------------------
# Merge the "Form" header cell across O2:O4 (was O2:O3) and give
# O4 the same bold / centred formatting already used by O3 so the
# whole merged block looks consistent.
# ---------------------------------------------------------------
$ws.Range("O4").Font.Bold = $true
$ws.Range("O4").HorizontalAlignment = -4108
$ws.Range("O4").IndentLevel = 0
$ws.Range("O2:O4").Merge()

# ---------------------------------------------------------------
# Column widths (values are expressed in the "ColumnWidth" COM
# units; Excel re-quantises these internally, so the numbers below
# are chosen to land as closely as possible on the widths used in
# the final workbook).
# ---------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 2.4986979166666665    # A
$ws.Columns.Item(2).ColumnWidth = 9.166666666666666     # B
$ws.Columns.Item(3).ColumnWidth = 3.9440104166666665    # C
$ws.Columns.Item(4).ColumnWidth = 24.053385416666668    # D
$ws.Range("E1:N1").EntireColumn.ColumnWidth = 9.166666666666666   # E..N
$ws.Columns.Item(15).ColumnWidth = 9.166666666666666    # O
$ws.Columns.Item(16).ColumnWidth = 34.276041666666664   # P

# ---------------------------------------------------------------
# Row heights
# ---------------------------------------------------------------
$ws.Rows.Item(2).RowHeight = 11.4
$ws.Rows.Item(3).RowHeight = 20.4
$ws.Rows.Item(4).RowHeight = 12
$ws.Rows.Item(6).RowHeight = 10.2
$ws.Rows.Item(7).RowHeight = 10.2

# ---------------------------------------------------------------
# Print area: A1:O7 -> A1:P7 (keep Print_Titles untouched)
# ---------------------------------------------------------------
foreach ($n in $wb.Names) {
    if ($n.Name -like "*Print_Area*") {
        $n.RefersTo = "=Sszw_Beispielbefüllung!`$A`$1:`$P`$7"
    }
}

# ---------------------------------------------------------------
# Selection / active cell (cosmetic view state)
# ---------------------------------------------------------------
$ws.Range("J23").Select()

Write-Host "sszw_vorlage update applied"
